$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 474.375
$ws.Range("J17").Value = 481.41025
$ws.Range("L17").Value = 1444.23075
$ws.Range("N17").Value = -1780.23075
# Row 33
$ws.Range("H33").Value = 536.02856
$ws.Range("I33").Value = 461.14816
$ws.Range("J33").Value = 788.75
$ws.Range("K33").Value = 461.14816
$ws.Range("L33").Value = 788.75
$ws.Range("M33").Value = -232.14816
$ws.Range("N33").Value = -1246.75
# Row 43
$ws.Range("H43").Value = 694.0833
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 694.0833
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 694.0833
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -832.0833
# Row 137
$ws.Range("H137").Value = 1314.2354
$ws.Range("I137").Value = 1127.4546
$ws.Range("K137").Value = 3382.3638
$ws.Range("M137").Value = -832.3638000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 32601.334
$ws.Range("J35").Value = 32601.334
$ws.Range("L35").Value = 32601.334
$ws.Range("N35").Value = -33221.334
# Row 64
$ws.Range("H64").Value = 385.30435
$ws.Range("I64").Value = 323.0909
$ws.Range("J64").Value = 442.33334
$ws.Range("K64").Value = 323.0909
$ws.Range("L64").Value = 442.33334
$ws.Range("M64").Value = -98.09089999999998
$ws.Range("N64").Value = -892.33334
# Row 67
$ws.Range("H67").Value = 385.30435
$ws.Range("I67").Value = 323.0909
$ws.Range("J67").Value = 442.33334
$ws.Range("K67").Value = 323.0909
$ws.Range("L67").Value = 442.33334
$ws.Range("M67").Value = 456.9091
$ws.Range("N67").Value = -2002.33334
# Row 80
$ws.Range("H80").Value = 512.5
$ws.Range("I80").Value = 652.8570999999999
$ws.Range("J80").Value = 423.18182
$ws.Range("K80").Value = 652.8570999999999
$ws.Range("L80").Value = 423.18182
$ws.Range("M80").Value = 345.1429000000001
$ws.Range("N80").Value = -2419.18182
# Row 82
$ws.Range("H82").Value = 18393.785
$ws.Range("J82").Value = 36209.332
$ws.Range("L82").Value = 36209.332
$ws.Range("N82").Value = -36975.332
# Row 83
$ws.Range("H83").Value = 512.5
$ws.Range("I83").Value = 652.8570999999999
$ws.Range("J83").Value = 423.18182
$ws.Range("K83").Value = 3264.2855
$ws.Range("L83").Value = 2115.9091
$ws.Range("M83").Value = 1727.7145
$ws.Range("N83").Value = -12099.9091
# Row 85
$ws.Range("H85").Value = 18393.785
$ws.Range("J85").Value = 36209.332
$ws.Range("L85").Value = 36209.332
$ws.Range("N85").Value = -38861.332
# Row 109
$ws.Range("H109").Value = 9690
$ws.Range("J109").Value = 9690
$ws.Range("L109").Value = 9690
$ws.Range("N109").Value = -12464

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1788.8889
$ws.Range("I16").Value = 1871.4286
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1871.4286
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1584.4286
$ws.Range("N16").Value = -2074
# Row 31
$ws.Range("H31").Value = 2301.2646
$ws.Range("I31").Value = 1792.0769
$ws.Range("J31").Value = 2616.476
$ws.Range("K31").Value = 1792.0769
$ws.Range("L31").Value = 2616.476
$ws.Range("M31").Value = -1497.0769
$ws.Range("N31").Value = -3206.476
# Row 34
$ws.Range("H34").Value = 2301.2646
$ws.Range("I34").Value = 1792.0769
$ws.Range("J34").Value = 2616.476
$ws.Range("K34").Value = 1792.0769
$ws.Range("L34").Value = 2616.476
$ws.Range("M34").Value = -1590.0769
$ws.Range("N34").Value = -3020.476
# Row 113
$ws.Range("H113").Value = 1788.8889
$ws.Range("I113").Value = 1871.4286
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1871.4286
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 298.5714
$ws.Range("N113").Value = -5840

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 590.73914
$ws.Range("I12").Value = 93.22221999999999
$ws.Range("J12").Value = 910.5714
$ws.Range("K12").Value = 279.66666
$ws.Range("L12").Value = 2731.7142
$ws.Range("M12").Value = -106.66666
$ws.Range("N12").Value = -3077.7142
# Row 68
$ws.Range("H68").Value = 1613.5955
$ws.Range("I68").Value = 1149.9524
$ws.Range("J68").Value = 2027.9149
$ws.Range("K68").Value = 3449.857199999999
$ws.Range("L68").Value = 6083.7447
$ws.Range("M68").Value = -2638.857199999999
$ws.Range("N68").Value = -7705.7447
# Row 71
$ws.Range("H71").Value = 1613.5955
$ws.Range("I71").Value = 1149.9524
$ws.Range("J71").Value = 2027.9149
$ws.Range("K71").Value = 10349.5716
$ws.Range("L71").Value = 18251.2341
$ws.Range("M71").Value = -6293.571599999999
$ws.Range("N71").Value = -26363.2341
# Row 113
$ws.Range("H113").Value = 475.51614
$ws.Range("I113").Value = 469.5
$ws.Range("J113").Value = 483.84616
$ws.Range("K113").Value = 1408.5
$ws.Range("L113").Value = 1451.53848
$ws.Range("M113").Value = 761.5
$ws.Range("N113").Value = -5791.53848

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 962.4091
$ws.Range("J107").Value = 1261.375
$ws.Range("L107").Value = 1261.375
$ws.Range("N107").Value = -5101.375
# Row 132
$ws.Range("H132").Value = 2423.111
$ws.Range("I132").Value = 1745.4814
$ws.Range("J132").Value = 3439.5557
$ws.Range("K132").Value = 5236.4442
$ws.Range("L132").Value = 10318.6671
$ws.Range("M132").Value = -2706.4442
$ws.Range("N132").Value = -15378.6671

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 62
$ws.Range("H62").Value = 34749.668
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 34749.668
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 34749.668
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -35997.668
# Row 65
$ws.Range("H65").Value = 34749.668
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 34749.668
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 104249.004
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -110489.004
# Row 94
$ws.Range("H94").Value = 21663.334
$ws.Range("J94").Value = 21663.334
$ws.Range("L94").Value = 21663.334
$ws.Range("N94").Value = -23015.334
# Row 127
$ws.Range("H127").Value = 47392
$ws.Range("J127").Value = 47392
$ws.Range("L127").Value = 47392
$ws.Range("N127").Value = -57312
# Row 132
$ws.Range("H132").Value = 1815.7941
$ws.Range("I132").Value = 1361.95
$ws.Range("J132").Value = 2464.1428
$ws.Range("K132").Value = 4085.85
$ws.Range("L132").Value = 7392.428400000001
$ws.Range("M132").Value = -1555.85
$ws.Range("N132").Value = -12452.4284

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1460.8182
$ws.Range("I107").Value = 981.9
$ws.Range("J107").Value = 6250
$ws.Range("K107").Value = 2945.7
$ws.Range("L107").Value = 18750
$ws.Range("M107").Value = -1025.7
$ws.Range("N107").Value = -22590
# Row 122
$ws.Range("H122").Value = 12999.19
$ws.Range("I122").Value = 20307.834
$ws.Range("J122").Value = 3254.3333
$ws.Range("K122").Value = 60923.50199999999
$ws.Range("L122").Value = 9762.999899999999
$ws.Range("M122").Value = -58473.50199999999
$ws.Range("N122").Value = -14662.9999
# Row 136
$ws.Range("H136").Value = 1515.2632
$ws.Range("I136").Value = 1169.697
$ws.Range("J136").Value = 3796
$ws.Range("K136").Value = 3509.090999999999
$ws.Range("L136").Value = 11388
$ws.Range("M136").Value = -959.0909999999994
$ws.Range("N136").Value = -16488
